# Rotated several components to fit manufacturer's placement
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

# Update rotation values (column E) for the affected components
$ws.Range("E38").Value = -90
$ws.Range("E39").Value = -90
$ws.Range("E57").Value = 90
$ws.Range("E59").Value = 180
$ws.Range("E60").Value = 180

# Reflect the updated view/scroll position and selection
$win.ScrollRow = 31
$win.ScrollColumn = 1
$ws.Range("E61").Select()
